$wb = $excel.ActiveWorkbook

$ws_Summary = $wb.Sheets.Item("Summary")
$ws_Summary.Range("B3").Value = 0.01
$ws_Summary.Range("B6").Value = 196497.6586768614
$ws_Summary.Range("B7").Value = 11121133.31629766
$ws_Summary.Range("B8").Value = 21897431.67196571
$ws_Summary.Range("B10").Value = 3788061.030312993

$ws_Fed_in_Capacity = $wb.Sheets.Item("Fed-in Capacity")
$ws_Fed_in_Capacity.Range("J5").Value = 107.2955742555736
$ws_Fed_in_Capacity.Range("K5").Value = 109.5572237694796
$ws_Fed_in_Capacity.Range("L5").Value = 98.64091687123928
$ws_Fed_in_Capacity.Range("M5").Value = 77.7676953375541
$ws_Fed_in_Capacity.Range("N5").Value = 74.3656454478664
$ws_Fed_in_Capacity.Range("O5").Value = 83.69133109099639
$ws_Fed_in_Capacity.Range("P5").Value = 106.2781106359148
$ws_Fed_in_Capacity.Range("Q5").Value = 128.4697750236904
$ws_Fed_in_Capacity.Range("J6").Value = 82.91243519753434
$ws_Fed_in_Capacity.Range("K6").Value = 62.76629045205057
$ws_Fed_in_Capacity.Range("L6").Value = 37.60657774285653
$ws_Fed_in_Capacity.Range("M6").Value = 24.33271034503603
$ws_Fed_in_Capacity.Range("N6").Value = 10.42253800004659
$ws_Fed_in_Capacity.Range("O6").Value = 31.97882363640291
$ws_Fed_in_Capacity.Range("P6").Value = 45.19417501179163
$ws_Fed_in_Capacity.Range("Q6").Value = 80.63453985745144
$ws_Fed_in_Capacity.Range("K7").Value = 83.72822537421013
$ws_Fed_in_Capacity.Range("L7").Value = 76.93542539304551
$ws_Fed_in_Capacity.Range("M7").Value = 77.82642397052864
$ws_Fed_in_Capacity.Range("N7").Value = 68.03899070462725
$ws_Fed_in_Capacity.Range("O7").Value = 83.36329197944329
$ws_Fed_in_Capacity.Range("P7").Value = 90.5862140395771
$ws_Fed_in_Capacity.Range("J8").Value = 92.06191329280011
$ws_Fed_in_Capacity.Range("K8").Value = 86.72591490550931
$ws_Fed_in_Capacity.Range("L8").Value = 70.31665934305323
$ws_Fed_in_Capacity.Range("M8").Value = 46.25150133729363
$ws_Fed_in_Capacity.Range("N8").Value = 42.33948647951536
$ws_Fed_in_Capacity.Range("O8").Value = 53.44993733213093
$ws_Fed_in_Capacity.Range("P8").Value = 80.46778072209531
$ws_Fed_in_Capacity.Range("Q8").Value = 109.0872921559756
$ws_Fed_in_Capacity.Range("J9").Value = 73.83937108913167
$ws_Fed_in_Capacity.Range("K9").Value = 47.25897875440141
$ws_Fed_in_Capacity.Range("L9").Value = 16.75508344424118
$ws_Fed_in_Capacity.Range("M9").Value = 0
$ws_Fed_in_Capacity.Range("O9").Value = 9.130000058643361
$ws_Fed_in_Capacity.Range("P9").Value = 26.85597970277001
$ws_Fed_in_Capacity.Range("Q9").Value = 68.37594193669796
$ws_Fed_in_Capacity.Range("K10").Value = 74.37427626277801
$ws_Fed_in_Capacity.Range("L10").Value = 64.96559098237566
$ws_Fed_in_Capacity.Range("M10").Value = 65.20591187910269
$ws_Fed_in_Capacity.Range("N10").Value = 55.71856618250388
$ws_Fed_in_Capacity.Range("O10").Value = 71.98338581476828
$ws_Fed_in_Capacity.Range("P10").Value = 80.84873837615825
$ws_Fed_in_Capacity.Range("L14").Value = 28.8362588753597
$ws_Fed_in_Capacity.Range("M14").Value = 0.09656339947136416
$ws_Fed_in_Capacity.Range("L16").Value = 47.4359690887577
$ws_Fed_in_Capacity.Range("N16").Value = 37.67551094615366
$ws_Fed_in_Capacity.Range("L35").Value = 28.83625887535973
$ws_Fed_in_Capacity.Range("M35").Value = 0.09656339947139259
$ws_Fed_in_Capacity.Range("L37").Value = 47.43596908875769
$ws_Fed_in_Capacity.Range("N37").Value = 37.67551094615368

$ws_Unmet_Demand = $wb.Sheets.Item("Unmet Demand")
$ws_Unmet_Demand.Range("G5").Value = 414.4337959369544
$ws_Unmet_Demand.Range("H5").Value = 330.5757541782243
$ws_Unmet_Demand.Range("I5").Value = 176.9760193775952
$ws_Unmet_Demand.Range("R5").Value = 95.28546653075742
$ws_Unmet_Demand.Range("S5").Value = 189.2190633734531
$ws_Unmet_Demand.Range("T5").Value = 219.2920578056454
$ws_Unmet_Demand.Range("U5").Value = 251.276137581582
$ws_Unmet_Demand.Range("G6").Value = 136.8785924310737
$ws_Unmet_Demand.Range("H6").Value = 107.7452501129632
$ws_Unmet_Demand.Range("I6").Value = 73.38935588968215
$ws_Unmet_Demand.Range("R6").Value = 71.29171788891175
$ws_Unmet_Demand.Range("S6").Value = 163.0473981187501
$ws_Unmet_Demand.Range("T6").Value = 198.2907557613397
$ws_Unmet_Demand.Range("U6").Value = 225.9107949275447
$ws_Unmet_Demand.Range("G7").Value = 167.6012020808691
$ws_Unmet_Demand.Range("H7").Value = 158.7616981666879
$ws_Unmet_Demand.Range("I7").Value = 143.7288091611985
$ws_Unmet_Demand.Range("J7").Value = 65.80192659108423
$ws_Unmet_Demand.Range("Q7").Value = 53.52351139825525
$ws_Unmet_Demand.Range("R7").Value = 159.7675876048201
$ws_Unmet_Demand.Range("S7").Value = 217.2238431175235
$ws_Unmet_Demand.Range("T7").Value = 226.2801774240348
$ws_Unmet_Demand.Range("U7").Value = 286.2977687777133
$ws_Unmet_Demand.Range("G8").Value = 414.2543098065221
$ws_Unmet_Demand.Range("H8").Value = 328.737591844935
$ws_Unmet_Demand.Range("I8").Value = 170.0563803341062
$ws_Unmet_Demand.Range("R8").Value = 84.01082088999308
$ws_Unmet_Demand.Range("S8").Value = 185.1290231762283
$ws_Unmet_Demand.Range("T8").Value = 218.5063572696783
$ws_Unmet_Demand.Range("U8").Value = 251.2617786911475
$ws_Unmet_Demand.Range("G9").Value = 136.7825588849655
$ws_Unmet_Demand.Range("H9").Value = 106.817768233445
$ws_Unmet_Demand.Range("I9").Value = 70.08293774516886
$ws_Unmet_Demand.Range("R9").Value = 65.32921403493339
$ws_Unmet_Demand.Range("S9").Value = 161.2636171197674
$ws_Unmet_Demand.Range("T9").Value = 197.9036731785617
$ws_Unmet_Demand.Range("U9").Value = 225.9044769310903
$ws_Unmet_Demand.Range("G10").Value = 167.5206907817879
$ws_Unmet_Demand.Range("H10").Value = 158.0458795257661
$ws_Unmet_Demand.Range("I10").Value = 141.3076148215571
$ws_Unmet_Demand.Range("J10").Value = 60.10977774604451
$ws_Unmet_Demand.Range("Q10").Value = 46.78178798155719
$ws_Unmet_Demand.Range("R10").Value = 156.1475068297698
$ws_Unmet_Demand.Range("S10").Value = 215.8207507508087
$ws_Unmet_Demand.Range("T10").Value = 225.9361746006879
$ws_Unmet_Demand.Range("U10").Value = 286.2933772523089

$ws_Household_Surplus = $wb.Sheets.Item("Household Surplus")
$ws_Household_Surplus.Range("B3").Value = 337587.5880544489
$ws_Household_Surplus.Range("B4").Value = 330152.41939151

$ws_Costs_and_Revenues = $wb.Sheets.Item("Costs and Revenues")
$ws_Costs_and_Revenues.Range("C2").Value = 116034.2185408279
$ws_Costs_and_Revenues.Range("D2").Value = 117549.91093637
$ws_Costs_and_Revenues.Range("C3").Value = 196825.9098199031
$ws_Costs_and_Revenues.Range("D3").Value = 38236.46568336456
$ws_Costs_and_Revenues.Range("E3").Value = 52530.53686621619
$ws_Costs_and_Revenues.Range("C4").Value = 47502.04170237896
$ws_Costs_and_Revenues.Range("D4").Value = 36277.64205058495
$ws_Costs_and_Revenues.Range("F4").Value = 21804.89414236139
$ws_Costs_and_Revenues.Range("M4").Value = 21804.8941423614
$ws_Costs_and_Revenues.Range("C5").Value = 38339.65294307929
$ws_Costs_and_Revenues.Range("D5").Value = 39312.96135688073
$ws_Costs_and_Revenues.Range("B6").Value = -43402.9036934624
$ws_Costs_and_Revenues.Range("C6").Value = -181357.7192440992
$ws_Costs_and_Revenues.Range("D6").Value = -10925.70685424888
$ws_Costs_and_Revenues.Range("E6").Value = 23785.86742824304
$ws_Costs_and_Revenues.Range("F6").Value = 76316.40429445921
$ws_Costs_and_Revenues.Range("G6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("H6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("I6").Value = 76316.40429445927
$ws_Costs_and_Revenues.Range("J6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("K6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("L6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("M6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("N6").Value = 76316.40429445924
$ws_Costs_and_Revenues.Range("O6").Value = 76316.40429445927
$ws_Costs_and_Revenues.Range("P6").Value = 76316.40429445924

$ws_Installed_Capacities = $wb.Sheets.Item("Installed Capacities")
$ws_Installed_Capacities.Range("C3").Value = 216.1492175724446
$ws_Installed_Capacities.Range("D3").Value = 260.7963925174648

$ws_Added_Capacities = $wb.Sheets.Item("Added Capacities")
$ws_Added_Capacities.Range("C3").Value = 216.1492175724445
$ws_Added_Capacities.Range("D3").Value = 44.64717494502023
$ws_Added_Capacities.Range("E3").Value = 65.38503947111997

$ws_PV_Dispatch = $wb.Sheets.Item("PV Dispatch")
$ws_PV_Dispatch.Range("G5").Value = 0.8689415781806812
$ws_PV_Dispatch.Range("H5").Value = 8.899047937542903
$ws_PV_Dispatch.Range("I5").Value = 33.49987019281074
$ws_PV_Dispatch.Range("J5").Value = 73.75033027111266
$ws_PV_Dispatch.Range("K5").Value = 110.5326272755009
$ws_PV_Dispatch.Range("L5").Value = 137.125498098748
$ws_PV_Dispatch.Range("M5").Value = 152.5785378897186
$ws_PV_Dispatch.Range("N5").Value = 155.0474181487245
$ws_PV_Dispatch.Range("O5").Value = 146.4068803306903
$ws_PV_Dispatch.Range("P5").Value = 124.9548851193548
$ws_PV_Dispatch.Range("Q5").Value = 93.8359148507591
$ws_PV_Dispatch.Range("R5").Value = 54.58365141039226
$ws_PV_Dispatch.Range("S5").Value = 19.80100621279229
$ws_PV_Dispatch.Range("T5").Value = 3.803791758485934
$ws_PV_Dispatch.Range("U5").Value = 0.06951532625445447
$ws_PV_Dispatch.Range("G6").Value = 0.4649247321369563
$ws_PV_Dispatch.Range("H6").Value = 4.490194123533237
$ws_PV_Dispatch.Range("I6").Value = 16.00727696173293
$ws_PV_Dispatch.Range("J6").Value = 43.92519146913236
$ws_PV_Dispatch.Range("K6").Value = 75.07514852230842
$ws_PV_Dispatch.Range("L6").Value = 100.9478020370177
$ws_PV_Dispatch.Range("M6").Value = 117.8013235769823
$ws_PV_Dispatch.Range("N6").Value = 120.9191740832867
$ws_PV_Dispatch.Range("O6").Value = 110.6174208080415
$ws_PV_Dispatch.Range("P6").Value = 88.78023240253862
$ws_PV_Dispatch.Range("Q6").Value = 59.34723422857008
$ws_PV_Dispatch.Range("R6").Value = 28.86611626373139
$ws_PV_Dispatch.Range("S6").Value = 8.635772985087758
$ws_PV_Dispatch.Range("T6").Value = 1.873972933481854
$ws_PV_Dispatch.Range("U6").Value = 0.03058715343006293
$ws_PV_Dispatch.Range("G7").Value = 0.3897772775896541
$ws_PV_Dispatch.Range("H7").Value = 3.465474340751655
$ws_PV_Dispatch.Range("I7").Value = 11.72166576605978
$ws_PV_Dispatch.Range("J7").Value = 27.55725352558855
$ws_PV_Dispatch.Range("K7").Value = 45.28503279632526
$ws_PV_Dispatch.Range("L7").Value = 57.94925088819277
$ws_PV_Dispatch.Range("M7").Value = 61.09935997707642
$ws_PV_Dispatch.Range("N7").Value = 59.64655376060593
$ws_PV_Dispatch.Range("O7").Value = 55.0932464723995
$ws_PV_Dispatch.Range("P7").Value = 47.14179000957051
$ws_PV_Dispatch.Range("Q7").Value = 32.63853185343913
$ws_PV_Dispatch.Range("R7").Value = 17.52580377234936
$ws_PV_Dispatch.Range("S7").Value = 6.792754919448789
$ws_PV_Dispatch.Range("T7").Value = 1.665412004246704
$ws_PV_Dispatch.Range("U7").Value = 0.02126057877761752
$ws_PV_Dispatch.Range("G8").Value = 1.048427708612923
$ws_PV_Dispatch.Range("H8").Value = 10.7372102708321
$ws_PV_Dispatch.Range("I8").Value = 40.41950923629976
$ws_PV_Dispatch.Range("J8").Value = 88.98399123388617
$ws_PV_Dispatch.Range("K8").Value = 133.3639361394712
$ws_PV_Dispatch.Range("L8").Value = 165.449755626934
$ws_PV_Dispatch.Range("M8").Value = 184.0947318899791
$ws_PV_Dispatch.Range("N8").Value = 187.0735771170756
$ws_PV_Dispatch.Range("O8").Value = 176.6482740895558
$ws_PV_Dispatch.Range("P8").Value = 150.7652150331742
$ws_PV_Dispatch.Range("Q8").Value = 113.2183977184739
$ws_PV_Dispatch.Range("R8").Value = 65.8582970511566
$ws_PV_Dispatch.Range("S8").Value = 23.89104641001701
$ws_PV_Dispatch.Range("T8").Value = 4.589492294453073
$ws_PV_Dispatch.Range("U8").Value = 0.08387421668903385
$ws_PV_Dispatch.Range("G9").Value = 0.560958278245113
$ws_PV_Dispatch.Range("H9").Value = 5.417676003051488
$ws_PV_Dispatch.Range("I9").Value = 19.31369510624622
$ws_PV_Dispatch.Range("J9").Value = 52.99825557753501
$ws_PV_Dispatch.Range("K9").Value = 90.58246021995758
$ws_PV_Dispatch.Range("L9").Value = 121.799296335633
$ws_PV_Dispatch.Range("M9").Value = 142.1340339220183
$ws_PV_Dispatch.Range("O9").Value = 133.4662443858011
$ws_PV_Dispatch.Range("P9").Value = 107.1184277115602
$ws_PV_Dispatch.Range("Q9").Value = 71.60583214932356
$ws_PV_Dispatch.Range("R9").Value = 34.82862011770975
$ws_PV_Dispatch.Range("S9").Value = 10.41955398407041
$ws_PV_Dispatch.Range("T9").Value = 2.261055516259907
$ws_PV_Dispatch.Range("U9").Value = 0.03690514988454693
$ws_PV_Dispatch.Range("G10").Value = 0.4702885766708382
$ws_PV_Dispatch.Range("H10").Value = 4.181292981673455
$ws_PV_Dispatch.Range("I10").Value = 14.14286010570121
$ws_PV_Dispatch.Range("J10").Value = 33.24940237062826
$ws_PV_Dispatch.Range("K10").Value = 54.63898190775737
$ws_PV_Dispatch.Range("L10").Value = 69.91908529886263
$ws_PV_Dispatch.Range("M10").Value = 73.71987206850237
$ws_PV_Dispatch.Range("N10").Value = 71.96697828272931
$ws_PV_Dispatch.Range("O10").Value = 66.47315263707451
$ws_PV_Dispatch.Range("P10").Value = 56.87926567298936
$ws_PV_Dispatch.Range("Q10").Value = 39.38025527013719
$ws_PV_Dispatch.Range("R10").Value = 21.14588454739968
$ws_PV_Dispatch.Range("S10").Value = 8.195847286163604
$ws_PV_Dispatch.Range("T10").Value = 2.009414827593581
$ws_PV_Dispatch.Range("U10").Value = 0.02565210418204575
$ws_PV_Dispatch.Range("M14").Value = 230.2496698278014
$ws_PV_Dispatch.Range("J16").Value = 41.58545896024957
$ws_PV_Dispatch.Range("L16").Value = 87.44870719248058
$ws_PV_Dispatch.Range("N16").Value = 90.01003351907953
$ws_PV_Dispatch.Range("M35").Value = 230.2496698278013
$ws_PV_Dispatch.Range("J37").Value = 41.58545896024956
$ws_PV_Dispatch.Range("L37").Value = 87.4487071924806
$ws_PV_Dispatch.Range("N37").Value = 90.01003351907951
